# Generate Report for Handoff
# Refresh the localization-status report with the new run's GUID-named
# files and updated timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "03c5bd04-b577-40e8-b656-1a019ea11ede"
$newGuid = "d963e7ae-6d89-4111-ae40-56042b4814f9"

$oldZhXlf = "$oldGuid.06b5930bf7aca12dde5fb9519a46ae3ee65b78b9.zh-cn.xlf"
$newZhXlf = "$newGuid.9ce56e505796f4f2d4779e1062d4ee92361e4367.zh-cn.xlf"

$oldDeXlf = "$oldGuid.06b5930bf7aca12dde5fb9519a46ae3ee65b78b9.de-de.xlf"
$newDeXlf = "$newGuid.9ce56e505796f4f2d4779e1062d4ee92361e4367.de-de.xlf"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cac3b57687f7243129f3007fdef1ba3997efba7b/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, "", "", "e2e\$newGuid.md")

$wsOverview.Range("G2").Value = "2016-08-16 10:57:56"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md")

$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-08-16 10:57:51"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md")

$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = "2016-08-16 10:57:56"
